# Apply updated cryptocurrency price/volume data scraped on
# Wed Dec 6 18:26:49 UTC 2023 to the "cryptos" worksheet (Sheet1).
# Several rows also changed rank/order (coin name + link swapped
# with the following row), which is reproduced by overwriting the
# Coin / Link / Price / Volume(1h) cell values directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values are plain numeric-looking strings (e.g. "231.70").
# Force those specific cells to Text format first so Excel keeps them
# as literal text (preserving trailing zeros) instead of coercing them
# into numbers.
$textCells = @('D5', 'D6', 'D7', 'D9', 'D10', 'D11', 'D12', 'D15', 'D16', 'D17', 'D20', 'D21', 'D22', 'D23', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D34', 'D35', 'D36', 'D37', 'D40', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '43.949.09'
$ws.Range('E2').Value = '  +0.58%  '
$ws.Range('D3').Value = '2.274.23'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '231.70'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('D6').Value = '0.626'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').Value = '63.63'
$ws.Range('E7').Value = '  +2.30%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '0.448'
$ws.Range('E9').Value = '  +7.01%  '
$ws.Range('D10').Value = '0.0995'
$ws.Range('E10').Value = '  +7.42%  '
$ws.Range('D11').Value = '57.58'
$ws.Range('E11').Value = '  -1.11%  '
$ws.Range('D12').Value = '27.25'
$ws.Range('E12').Value = '  +14.79%  '
$ws.Range('D14').Value = '2.613.15'
$ws.Range('E14').Value = '  -1.13%  '
$ws.Range('D15').Value = '15.73'
$ws.Range('E15').Value = '  -1.37%  '
$ws.Range('D16').Value = '6.13'
$ws.Range('E16').Value = '  +6.52%  '
$ws.Range('D17').Value = '0.835'
$ws.Range('E17').Value = '  +2.05%  '
$ws.Range('D18').Value = '2.274.35'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('D19').Value = '43.910.70'
$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('D20').Value = '0.0000100'
$ws.Range('E20').Value = '  +7.38%  '
$ws.Range('D21').Value = '73.81'
$ws.Range('E21').Value = '  +0.30%  '
$ws.Range('D22').Value = '6.10'
$ws.Range('E22').Value = '  -3.04%  '
$ws.Range('D23').Value = '252.80'
$ws.Range('E23').Value = '  +1.01%  '
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('E25').Value = '  -4.83%  '
$ws.Range('B26').Value = 'WEMIXToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D26').Value = '3.35'
$ws.Range('E26').Value = '  +25.35%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').Value = '2.26'
$ws.Range('E27').Value = '  -5.27%  '
$ws.Range('D28').Value = '10.06'
$ws.Range('E28').Value = '  +1.91%  '
$ws.Range('D29').Value = '171.71'
$ws.Range('E29').Value = '  +0.69%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = '20.89'
$ws.Range('E30').Value = '  +1.21%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').Value = '0.138'
$ws.Range('E31').Value = '  -3.38%  '
$ws.Range('E32').Value = '  -5.87%  '
$ws.Range('E33').Value = '  +2.04%  '
$ws.Range('D34').Value = '0.0698'
$ws.Range('E34').Value = '  +5.85%  '
$ws.Range('D35').Value = '4.81'
$ws.Range('E35').Value = '  +0.78%  '
$ws.Range('D36').Value = '4.89'
$ws.Range('E36').Value = '  -3.23%  '
$ws.Range('D37').Value = '3.81'
$ws.Range('E37').Value = '  +3.60%  '
$ws.Range('E38').Value = '  -0.63%  '
$ws.Range('E39').Value = '  -5.96%  '
$ws.Range('D40').Value = '0.0257'
$ws.Range('E40').Value = '  +2.50%  '
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('D42').Value = '0.000229'
$ws.Range('E42').Value = '  +5.56%  '
$ws.Range('D43').Value = '17.70'
$ws.Range('E43').Value = '  +4.35%  '
$ws.Range('D44').Value = '0.0975'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '8.27'
$ws.Range('E45').Value = '  -6.36%  '
$ws.Range('B46').Value = 'Celestia'
$ws.Range('C46').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D46').Value = '10.44'
$ws.Range('E46').Value = '  +13.28%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').Value = '1.21'
$ws.Range('E47').Value = '  -0.61%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '98.17'
$ws.Range('E48').Value = '  -0.36%  '
$ws.Range('B49').Value = 'FTXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D49').Value = '4.40'
$ws.Range('E49').Value = '  -5.64%  '
$ws.Range('D50').Value = '1.447.31'
$ws.Range('E50').Value = '  -1.81%  '
$ws.Range('D51').Value = '2.31'
$ws.Range('E51').Value = '  +0.47%  '
